$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "priceBF" field description (row 47) to describe Betfair Starting Price.
$ws.Range("B47").Value = "Betfair Starting Price (not available before the race)"

# Update the "priceSP" field description (row 45) to mention Pre-Post price.
$ws.Range("B45").Value = "Official starting price (or Pre-Post price if before the race)"

# Move the selection/active cell to B46 and scroll the window so row 38 is at
# the top, matching the saved view state of the edited workbook.
$ws.Range("B46").Select()
$excel.ActiveWindow.ScrollRow = 38
